$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.647.49'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.77%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.943.78'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '483.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.54'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.727'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.66%  '
$ws.Range('E10').Value = '  +8.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000354'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +11.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.66'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.96%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.55'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.569.37'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.995.31'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.71'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.82'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('E19').Value = '  -2.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.647.93'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '434.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.55%  '
$ws.Range('E23').Value = '  -2.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.27'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.48%  '
$ws.Range('E26').Value = '  +1.79%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '38.36'
$ws.Range('D27').Style = 'Normal'
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('E29').Value = '  +6.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '719.80'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.29'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.01%  '
$ws.Range('E32').Value = '  -5.05%  '
$ws.Range('E33').Value = '  +3.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0894'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +31.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.85'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.18'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.151'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.58'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0470'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.50%  '
$ws.Range('E42').Value = '  +10.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.343'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.45%  '
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('B47').Value = 'LidoDAOToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.24'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '146.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('E51').Value = '  -1.27%  '
